$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

$t.Cell(1, 1).Range.Text = "0M"
$t.Cell(2, 1).Range.Text = "0M"
$t.Cell(3, 1).Range.Text = "0M"
$t.Cell(4, 1).Range.Text = "6087"

$t.Cell(6, 1).Range.Text = "0.02271"
$t.Cell(7, 1).Range.Text = "0.00361"
$t.Cell(8, 1).Range.Text = "0.00039"
$t.Cell(9, 1).Range.Text = "0.01659"
$t.Cell(10, 1).Range.Text = "0.01708"
$t.Cell(11, 1).Range.Text = "0.01744"
$t.Cell(12, 1).Range.Text = "1.29336"

$t.Cell(44, 1).Range.Text = "99.97"
$t.Cell(45, 1).Range.Text = "1.29"
$t.Cell(46, 1).Range.Text = "3840"
